$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 667.4761999999999
$ws.Cells.Item(28, 9).Value = 665.41174
$ws.Cells.Item(28, 10).Value = 676.25
$ws.Cells.Item(28, 11).Value = 665.41174
$ws.Cells.Item(28, 12).Value = 676.25
$ws.Cells.Item(28, 13).Value = -180.41174
$ws.Cells.Item(28, 14).Value = -1646.25
$ws.Cells.Item(62, 8).Value = 3166.2
$ws.Cells.Item(62, 9).Value = 2499.5454
$ws.Cells.Item(62, 10).Value = 4999.5
$ws.Cells.Item(62, 11).Value = 2499.5454
$ws.Cells.Item(62, 12).Value = 4999.5
$ws.Cells.Item(62, 13).Value = -1875.5454
$ws.Cells.Item(62, 14).Value = -6247.5
$ws.Cells.Item(65, 8).Value = 3166.2
$ws.Cells.Item(65, 9).Value = 2499.5454
$ws.Cells.Item(65, 10).Value = 4999.5
$ws.Cells.Item(65, 11).Value = 12497.727
$ws.Cells.Item(65, 12).Value = 24997.5
$ws.Cells.Item(65, 13).Value = -9377.726999999999
$ws.Cells.Item(65, 14).Value = -31237.5
$ws.Cells.Item(100, 8).Value = 40001616
$ws.Cells.Item(100, 9).Value = 40001616
$ws.Cells.Item(100, 11).Value = 40001616
$ws.Cells.Item(100, 13).Value = -40001075
$ws.Cells.Item(116, 8).Value = 509240.8
$ws.Cells.Item(116, 9).Value = 1112912.2
$ws.Cells.Item(116, 10).Value = 15327.818
$ws.Cells.Item(116, 11).Value = 1112912.2
$ws.Cells.Item(116, 12).Value = 15327.818
$ws.Cells.Item(116, 13).Value = -1109470.2
$ws.Cells.Item(116, 14).Value = -22211.818

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 539.5217
$ws.Cells.Item(2, 9).Value = 508.60526
$ws.Cells.Item(2, 10).Value = 686.375
$ws.Cells.Item(2, 11).Value = 508.60526
$ws.Cells.Item(2, 12).Value = 686.375
$ws.Cells.Item(2, 13).Value = -395.60526
$ws.Cells.Item(2, 14).Value = -912.375
$ws.Cells.Item(61, 8).Value = 1378.317
$ws.Cells.Item(61, 9).Value = 1167.5135
$ws.Cells.Item(61, 11).Value = 1167.5135
$ws.Cells.Item(61, 13).Value = -955.5135
$ws.Cells.Item(97, 8).Value = 767.56525
$ws.Cells.Item(97, 9).Value = 636.6667
$ws.Cells.Item(97, 10).Value = 1013
$ws.Cells.Item(97, 11).Value = 636.6667
$ws.Cells.Item(97, 12).Value = 1013
$ws.Cells.Item(97, 13).Value = -140.6667
$ws.Cells.Item(97, 14).Value = -2005
$ws.Cells.Item(110, 8).Value = 3564
$ws.Cells.Item(110, 9).Value = 4602.4
$ws.Cells.Item(110, 10).Value = 1833.3334
$ws.Cells.Item(110, 11).Value = 4602.4
$ws.Cells.Item(110, 12).Value = 1833.3334
$ws.Cells.Item(110, 13).Value = -2557.4
$ws.Cells.Item(110, 14).Value = -5923.3334
$ws.Cells.Item(116, 8).Value = 539.5217
$ws.Cells.Item(116, 9).Value = 508.60526
$ws.Cells.Item(116, 10).Value = 686.375
$ws.Cells.Item(116, 11).Value = 508.60526
$ws.Cells.Item(116, 12).Value = 686.375
$ws.Cells.Item(116, 13).Value = 1785.39474
$ws.Cells.Item(116, 14).Value = -5274.375
$ws.Cells.Item(136, 8).Value = 1378.317
$ws.Cells.Item(136, 9).Value = 1167.5135
$ws.Cells.Item(136, 11).Value = 3502.5405
$ws.Cells.Item(136, 13).Value = -952.5405000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 539.5217
$ws.Cells.Item(3, 9).Value = 508.60526
$ws.Cells.Item(3, 10).Value = 686.375
$ws.Cells.Item(3, 11).Value = 508.60526
$ws.Cells.Item(3, 12).Value = 686.375
$ws.Cells.Item(3, 13).Value = -394.60526
$ws.Cells.Item(3, 14).Value = -914.375
$ws.Cells.Item(129, 8).Value = 44699.715
$ws.Cells.Item(129, 10).Value = 44699.715
$ws.Cells.Item(129, 12).Value = 44699.715
$ws.Cells.Item(129, 14).Value = -54699.715
$ws.Cells.Item(134, 8).Value = 2407.836
$ws.Cells.Item(134, 9).Value = 1456.4386
$ws.Cells.Item(134, 10).Value = 7830.8
$ws.Cells.Item(134, 11).Value = 4369.3158
$ws.Cells.Item(134, 12).Value = 23492.4
$ws.Cells.Item(134, 13).Value = -1834.3158
$ws.Cells.Item(134, 14).Value = -28562.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(86, 8).Value = 2914.1428
$ws.Cells.Item(86, 9).Value = 1950
$ws.Cells.Item(86, 10).Value = 3299.8
$ws.Cells.Item(86, 11).Value = 1950
$ws.Cells.Item(86, 12).Value = 3299.8
$ws.Cells.Item(86, 13).Value = -827
$ws.Cells.Item(86, 14).Value = -5545.8
$ws.Cells.Item(89, 8).Value = 2914.1428
$ws.Cells.Item(89, 9).Value = 1950
$ws.Cells.Item(89, 10).Value = 3299.8
$ws.Cells.Item(89, 11).Value = 9750
$ws.Cells.Item(89, 12).Value = 16499
$ws.Cells.Item(89, 13).Value = -4134
$ws.Cells.Item(89, 14).Value = -27731
$ws.Cells.Item(122, 8).Value = 2499.56
$ws.Cells.Item(122, 10).Value = 4360
$ws.Cells.Item(122, 12).Value = 13080
$ws.Cells.Item(122, 14).Value = -17980

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(56, 8).Value = 5429.909
$ws.Cells.Item(56, 9).Value = 5429.909
$ws.Cells.Item(56, 11).Value = 5429.909
$ws.Cells.Item(56, 13).Value = -4899.909
$ws.Cells.Item(131, 8).Value = 10417638
$ws.Cells.Item(131, 10).Value = 820.3953
$ws.Cells.Item(131, 12).Value = 2461.1859
$ws.Cells.Item(131, 14).Value = -12541.1859

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(93, 8).Value = 19400
$ws.Cells.Item(93, 10).Value = 19400
$ws.Cells.Item(93, 12).Value = 19400
$ws.Cells.Item(93, 14).Value = -23144
$ws.Cells.Item(112, 8).Value = 28556.666
$ws.Cells.Item(112, 10).Value = 28556.666
$ws.Cells.Item(112, 12).Value = 28556.666
$ws.Cells.Item(112, 14).Value = -30772.666
$ws.Cells.Item(126, 8).Value = 3796.35
$ws.Cells.Item(126, 9).Value = 2741.5574
$ws.Cells.Item(126, 11).Value = 8224.672200000001
$ws.Cells.Item(126, 13).Value = -5754.672200000001
$ws.Cells.Item(136, 8).Value = 17081.5
$ws.Cells.Item(136, 10).Value = 17081.5
$ws.Cells.Item(136, 12).Value = 51244.5
$ws.Cells.Item(136, 14).Value = -56344.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6172.75
$ws.Cells.Item(7, 9).Value = 3799.6667
$ws.Cells.Item(7, 10).Value = 8545.833000000001
$ws.Cells.Item(7, 11).Value = 3799.6667
$ws.Cells.Item(7, 12).Value = 8545.833000000001
$ws.Cells.Item(7, 13).Value = -3687.6667
$ws.Cells.Item(7, 14).Value = -8769.833000000001
$ws.Cells.Item(111, 8).Value = 37400
$ws.Cells.Item(111, 10).Value = 37400
$ws.Cells.Item(111, 12).Value = 37400
$ws.Cells.Item(111, 14).Value = -45580
$ws.Cells.Item(126, 8).Value = 6172.75
$ws.Cells.Item(126, 9).Value = 3799.6667
$ws.Cells.Item(126, 10).Value = 8545.833000000001
$ws.Cells.Item(126, 11).Value = 11399.0001
$ws.Cells.Item(126, 12).Value = 25637.499
$ws.Cells.Item(126, 13).Value = -8929.000100000001
$ws.Cells.Item(126, 14).Value = -30577.499
$ws.Cells.Item(132, 8).Value = 4017.0518
$ws.Cells.Item(132, 9).Value = 1544.68
$ws.Cells.Item(132, 10).Value = 5890.0605
$ws.Cells.Item(132, 11).Value = 4634.04
$ws.Cells.Item(132, 12).Value = 17670.1815
$ws.Cells.Item(132, 13).Value = -2104.04
$ws.Cells.Item(132, 14).Value = -22730.1815
$ws.Cells.Item(136, 8).Value = 2774.1462
$ws.Cells.Item(136, 9).Value = 1475.5
$ws.Cells.Item(136, 10).Value = 4277.8423
$ws.Cells.Item(136, 11).Value = 4426.5
$ws.Cells.Item(136, 12).Value = 12833.5269
$ws.Cells.Item(136, 13).Value = -1876.5
$ws.Cells.Item(136, 14).Value = -17933.5269

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(108, 8).Value = 39750
$ws.Cells.Item(108, 10).Value = 39750
$ws.Cells.Item(108, 12).Value = 39750
$ws.Cells.Item(108, 14).Value = -47430
$ws.Cells.Item(126, 8).Value = 5569.683
$ws.Cells.Item(126, 9).Value = 2542.9033
$ws.Cells.Item(126, 10).Value = 14952.7
$ws.Cells.Item(126, 11).Value = 7628.7099
$ws.Cells.Item(126, 12).Value = 44858.10000000001
$ws.Cells.Item(126, 13).Value = -5158.7099
$ws.Cells.Item(126, 14).Value = -49798.10000000001
$ws.Cells.Item(136, 8).Value = 3480.3823
$ws.Cells.Item(136, 9).Value = 1069.7778
$ws.Cells.Item(136, 11).Value = 3209.3334
$ws.Cells.Item(136, 13).Value = -659.3334000000004
